$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the last data row (row 6); this also shrinks the used range to A1:AH5.
$ws.Rows.Item(6).Delete()

# 2) Replace the data values in rows 2-5 (columns A..AH) with the new dataset.
$data = @{
    2 = @{ 1=45035.50694444445; 2=5.237; 3=1.607; 4=0; 5=2.478; 6=3.053; 7=2.519; 8=5.488; 9=1.554; 10=0.941; 11=4.022; 12=1.069; 13=0.938; 14=0.693; 15=0.87; 16=2.749; 17=1.106; 18=0.51; 19=0.063; 20=20.094; 21=4.803; 22=2.372; 23=3.891; 24=1.034; 25=0.249; 26=1.98; 27=1.154; 28=0.674; 29=0.94; 30=3.129; 31=2.866; 32=3.128; 33=0.446; 34=1.43 }
    3 = @{ 1=45035.51388888889; 2=18.947; 3=13.45; 4=0.433; 5=37.331; 6=31.543; 7=14.468; 8=47.36; 9=20.827; 10=9.587; 11=15.327; 12=15.113; 13=15.715; 14=4.481; 15=13.404; 16=20.245; 17=11.26; 18=0.521; 19=0.431; 20=203.981; 21=38.714; 22=13.174; 23=26.433; 24=13.46; 25=1.824; 26=24.234; 27=11.25; 28=9.607; 29=11.408; 30=16.915; 31=1.133; 32=42.204; 33=7.027; 34=15.704 }
    4 = @{ 1=45035.52083333334; 2=15.631; 3=11.298; 4=0.357; 5=31.511; 6=26.53; 7=12.077; 8=46.435; 9=17.544; 10=8.078; 11=12.652; 12=12.776; 13=13.258; 14=3.75; 15=11.325; 16=16.994; 17=9.466; 18=0.414; 19=0.346; 20=169.829; 21=32.624; 22=10.981; 23=22.193; 24=11.356; 25=1.528; 26=22.662; 27=9.462; 28=8.108; 29=9.609; 30=14.087; 31=0.716; 32=41.951; 33=5.931; 34=13.222 }
    5 = @{ 1=45035.52777777778; 2=11.59; 3=8.4; 4=0.25; 5=23.37; 6=19.69; 7=8.98; 8=35.84; 9=13; 10=6; 11=9.37; 12=9.5; 13=9.82; 14=2.78; 15=8.41; 16=12.64; 17=7.02; 18=0.33; 19=0.24; 20=124.12; 21=24.24; 22=8.15; 23=16.47; 24=8.43; 25=1.13; 26=17.21; 27=7.03; 28=6.01; 29=7.13; 30=10.46; 31=0.52; 32=32.33; 33=4.4; 34=9.81 }
}

foreach ($r in $data.Keys) {
    $rowData = $data[$r]
    foreach ($c in $rowData.Keys) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c]
    }
}

# 3) Widen a subset of the data columns by one character.
#    ColumnWidth values below are chosen so the stored OOXML <col width="..">
#    lands exactly on the target integer width (7 -> 8, 8 -> 9).
$colsTo8 = @(2,3,7,9,11,12,13,15,16,22,23,24,26,29,30,34)
foreach ($c in $colsTo8) {
    $ws.Columns.Item($c).ColumnWidth = 7.165
}

$colsTo9 = @(20)
foreach ($c in $colsTo9) {
    $ws.Columns.Item($c).ColumnWidth = 8.165
}
